$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 485.9655
$ws.Range("J17").Value = 485.9655
$ws.Range("L17").Value = 1457.8965
$ws.Range("N17").Value = -1793.8965
$ws.Range("H33").Value = 314.5238
$ws.Range("I33").Value = 297.57144
$ws.Range("J33").Value = 348.42856
$ws.Range("K33").Value = 297.57144
$ws.Range("L33").Value = 348.42856
$ws.Range("M33").Value = -68.57144
$ws.Range("N33").Value = -806.4285600000001
$ws.Range("H98").Value = 781.2778
$ws.Range("I98").Value = 758.93335
$ws.Range("K98").Value = 758.93335
$ws.Range("M98").Value = 739.06665
$ws.Range("H122").Value = 781.2778
$ws.Range("I122").Value = 758.93335
$ws.Range("K122").Value = 2276.80005
$ws.Range("M122").Value = 173.1999500000002
$ws.Range("H138").Value = 2490297.2
$ws.Range("I138").Value = 1478.4584
$ws.Range("J138").Value = 3879405.5
$ws.Range("K138").Value = 4435.3752
$ws.Range("L138").Value = 11638216.5
$ws.Range("M138").Value = 704.6247999999996
$ws.Range("N138").Value = -11648496.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2041.82
$ws.Range("I32").Value = 1818.1818
$ws.Range("J32").Value = 3681.8333
$ws.Range("K32").Value = 1818.1818
$ws.Range("L32").Value = 3681.8333
$ws.Range("M32").Value = -1531.1818
$ws.Range("N32").Value = -4255.8333
$ws.Range("H45").Value = 2765
$ws.Range("I45").Value = 10000
$ws.Range("J45").Value = 1318
$ws.Range("K45").Value = 10000
$ws.Range("L45").Value = 1318
$ws.Range("M45").Value = -9623
$ws.Range("N45").Value = -2072
$ws.Range("H132").Value = 55406.566
$ws.Range("I132").Value = 33098.227
$ws.Range("J132").Value = 170666.33
$ws.Range("K132").Value = 99294.681
$ws.Range("L132").Value = 511998.99
$ws.Range("M132").Value = -96764.681
$ws.Range("N132").Value = -517058.99

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1381.6364
$ws.Range("I107").Value = 1496
$ws.Range("J107").Value = 992.8
$ws.Range("K107").Value = 1496
$ws.Range("L107").Value = 992.8
$ws.Range("M107").Value = 424
$ws.Range("N107").Value = -4832.8
$ws.Range("H134").Value = 2105.5083
$ws.Range("I134").Value = 1627.6731
$ws.Range("J134").Value = 4866.3335
$ws.Range("K134").Value = 4883.0193
$ws.Range("L134").Value = 14599.0005
$ws.Range("M134").Value = -2348.0193
$ws.Range("N134").Value = -19669.0005
$ws.Range("H137").Value = 50000
$ws.Range("J137").Value = 50000
$ws.Range("L137").Value = 50000
$ws.Range("N137").Value = -60200
$ws.Range("H138").Value = 32475
$ws.Range("J138").Value = 32475
$ws.Range("L138").Value = 32475
$ws.Range("N138").Value = -42755
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 2000
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 2000
$ws.Range("H58").Value = 23257334
$ws.Range("I58").Value = 25001330
$ws.Range("J58").Value = 4038.3333
$ws.Range("K58").Value = 25001330
$ws.Range("L58").Value = 4038.3333
$ws.Range("M58").Value = -25001127
$ws.Range("N58").Value = -4444.3333
$ws.Range("H99").Value = 2501.2122
$ws.Range("I99").Value = 2382.4
$ws.Range("J99").Value = 2600.2222
$ws.Range("K99").Value = 2382.4
$ws.Range("L99").Value = 2600.2222
$ws.Range("M99").Value = -884.4000000000001
$ws.Range("N99").Value = -5596.2222
$ws.Range("H105").Value = 765.95
$ws.Range("I105").Value = 836.13336
$ws.Range("J105").Value = 555.4
$ws.Range("K105").Value = 836.13336
$ws.Range("L105").Value = 555.4
$ws.Range("M105").Value = 910.86664
$ws.Range("N105").Value = -4049.4
$ws.Range("H107").Value = 352.72726
$ws.Range("I107").Value = 334.54544
$ws.Range("K107").Value = 334.54544
$ws.Range("M107").Value = 1585.45456
$ws.Range("H126").Value = 2501.2122
$ws.Range("I126").Value = 2382.4
$ws.Range("J126").Value = 2600.2222
$ws.Range("K126").Value = 7147.200000000001
$ws.Range("L126").Value = 7800.6666
$ws.Range("M126").Value = -4677.200000000001
$ws.Range("N126").Value = -12740.6666
$ws.Range("H132").Value = 39959.332
$ws.Range("I132").Value = 2634.5
$ws.Range("J132").Value = 146601.72
$ws.Range("K132").Value = 7903.5
$ws.Range("L132").Value = 439805.16
$ws.Range("M132").Value = -5373.5
$ws.Range("N132").Value = -444865.16
$ws.Range("H134").Value = 21662.697
$ws.Range("I134").Value = 1511.94
$ws.Range("J134").Value = 189585.67
$ws.Range("K134").Value = 4535.82
$ws.Range("L134").Value = 568757.01
$ws.Range("M134").Value = -2000.82
$ws.Range("N134").Value = -573827.01
$ws.Range("H136").Value = 23257334
$ws.Range("I136").Value = 25001330
$ws.Range("J136").Value = 4038.3333
$ws.Range("K136").Value = 75003990
$ws.Range("L136").Value = 12114.9999
$ws.Range("M136").Value = -75001440
$ws.Range("N136").Value = -17214.9999
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -2226

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 66667708
$ws.Range("I12").Value = 166668530
$ws.Range("J12").Value = 498.66666
$ws.Range("K12").Value = 500005590
$ws.Range("L12").Value = 1495.99998
$ws.Range("M12").Value = -500005417
$ws.Range("N12").Value = -1841.99998
$ws.Range("H87").Value = 25410.318
$ws.Range("I87").Value = 17002
$ws.Range("J87").Value = 29334.2
$ws.Range("K87").Value = 51006
$ws.Range("L87").Value = 88002.60000000001
$ws.Range("M87").Value = -49758
$ws.Range("N87").Value = -90498.60000000001
$ws.Range("H90").Value = 25410.318
$ws.Range("I90").Value = 17002
$ws.Range("J90").Value = 29334.2
$ws.Range("K90").Value = 153018
$ws.Range("L90").Value = 264007.8
$ws.Range("M90").Value = -146778
$ws.Range("N90").Value = -276487.8
$ws.Range("H131").Value = 982.3678
$ws.Range("I131").Value = 564.2222
$ws.Range("J131").Value = 1030.6154
$ws.Range("K131").Value = 1692.6666
$ws.Range("L131").Value = 3091.8462
$ws.Range("M131").Value = 3347.3334
$ws.Range("N131").Value = -13171.8462

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 887.7857
$ws.Range("I102").Value = 782.9
$ws.Range("J102").Value = 1150
$ws.Range("K102").Value = 782.9
$ws.Range("L102").Value = 1150
$ws.Range("M102").Value = 839.1
$ws.Range("N102").Value = -4394
$ws.Range("H132").Value = 57266.25
$ws.Range("I132").Value = 37110.25
$ws.Range("J132").Value = 127812.25
$ws.Range("K132").Value = 111330.75
$ws.Range("L132").Value = 383436.75
$ws.Range("M132").Value = -108800.75
$ws.Range("N132").Value = -388496.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 387.9091
$ws.Range("I55").Value = 207
$ws.Range("J55").Value = 455.75
$ws.Range("K55").Value = 207
$ws.Range("L55").Value = 455.75
$ws.Range("M55").Value = -34
$ws.Range("N55").Value = -801.75
$ws.Range("H122").Value = 3266.5908
$ws.Range("I122").Value = 2925
$ws.Range("J122").Value = 3461.7856
$ws.Range("K122").Value = 8775
$ws.Range("L122").Value = 10385.3568
$ws.Range("M122").Value = -6325
$ws.Range("N122").Value = -15285.3568
